$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 164 (shifts old rows 164..270 down to 165..271,
# carrying their formatting with them).
$ws.Rows(164).Insert()

# Populate the newly inserted row 164 with the new data record.
$ws.Range("A164").Value = 10
$ws.Range("B164").Value = "Vega Modelo de Temuco"
$ws.Range("C164").Value = "La Araucanía"
$ws.Range("D164").Value = 44603
$ws.Range("E164").Value = 9
$ws.Range("F164").Value = 100112044
$ws.Range("G164").Value = "Perejil"
$ws.Range("H164").Value = "Sin especificar"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 75
$ws.Range("K164").Value = 5000
$ws.Range("L164").Value = 6000
$ws.Range("M164").Value = 5533
$ws.Range("N164").Value = "$/docena de atados (3 kilos)"
$ws.Range("O164").Value = "Provincia de Cautín"
$ws.Range("P164").Value = 1844
$ws.Range("Q164").Value = 3
$ws.Range("R164").Value = "Hortaliza"
